# Updates the "hashcode" column (column B) for a set of rows identified by
# their code in column A. This mirrors an automated hashcode refresh where
# only the hash values change while the codes they belong to stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("hashcode.csv")

# Map of code (column A) -> new hashcode value (column B)
$updates = @{
    "05-050301A"           = "d2eb60f23630c81630eb9f1212449d63"
    "05-050207TP"          = "1cf29da4f485f4f8d6a2a0ab20e6fa0a"
    "05-050316TC"          = "83c5f6634eb2e2e8396ebdffa05479f5"
    "05-050301TP"          = "eebce30625319f1dfc293e782e87b5df"
    "05-050301TC"          = "d0abdc7b18468b17338332db5ae8a8ea"
    "05-050312TP"          = "8cfb26f044370eadd9dacd6571c0c16a"
    "05-050203TP"          = "47976b99ccdc02fbe30f20bc5d4c66e4"
    "05-050203TC"          = "4069a1af8c87d42f46a99d106e8a998b"
    "05-050314TP"          = "cdee2c75c757919e59b2193851b6f441"
    "05-050314TC"          = "ddf74f49b35251288f6e9e82cae9b7cb"
    "05-050205TP"          = "5a4da336cfaddcd1b3cc63fdbe620b3f"
    "05-050205TC"          = "2662d2bb210dba6c8a0adcd5d9fba2d7"
    "01-080101-010112TM"   = "0e4158b3be5756e9866cace2776c8db4"
    "05-050201TC"          = "711312510193327a95337f98f5a54eb4"
    "05-050201TP"          = "7e47f44727357eb02834ea30b6213b61"
    "05-050204A"           = "0ed3cebe67051283f8a2c674dcc51603"
    "05-050313A"           = "7ede5c43201f39fc1beae5e86f411e96"
    "05-050205A"           = "12e757b398212a1702dda98cca6a66ce"
    "05-050314A"           = "9f7edac59a821b3bc92b16054bd464ea"
    "05-050312A"           = "fbe1d95224c2b8044ace0593cde515a9"
    "05-050203A"           = "9b96c17f4564f2c119b882f8b42d6447"
    "05-050204TP"          = "2529bbe5f46938cd38c5557b2739e83a"
    "05-050204TC"          = "9df70dbea5d1a009415c0601313feacb"
    "05-050302TC"          = "932d7c111d75136c6a37a83e88d7afa6"
    "05-050313TP"          = "1543eba71c99b7491d6f82816b926040"
    "05-050206TP"          = "366d6593038ea112c8c73173a90f2cf6"
    "05-050206TC"          = "d4901b5ec07c3cc19aff9ec86f469438"
    "05-050206A"           = "ebf3f6b8d66231f35dda54cea2a650a4"
    "05-050315A"           = "79d91e5251c5085a3df550195b408243"
    "05-050207A"           = "942850043d43e3ef28b1e586dfc9ca8d"
    "05-050315TP"          = "280bc64f442fba14a08af0482de06fe3"
    "05-050202TC"          = "45364393e1315b47478d23e7ecbfdf31"
    "05-050311TC"          = "96a5c694513a1f5f4eba3b629595810f"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    if ($null -ne $code -and $updates.ContainsKey([string]$code)) {
        $ws.Cells.Item($r, 2).Value = $updates[[string]$code]
    }
}
